# ---------------------------------------------------------------------------
# Trust_Basic.docx – SharePoint custom-XML metadata refresh
#
# The source diff only touches the document's *custom XML* parts (the
# SharePoint "document management" metadata that Word/SharePoint keeps in
# customXml/item2.xml, itemProps2.xml, item3.xml, etc.):
#
#   * customXml/item2.xml   – the content-type schema part gets a new
#                              ma:versionID and a new ma:fieldsID (both are
#                              content hashes SharePoint mints whenever the
#                              associated field schema changes).
#   * customXml/itemProps2.xml – re-minted with a fresh ds:itemID and an
#                              emptied <ds:schemaRefs> list.
#   * customXml/item4.xml / itemProps4.xml – a brand new custom XML part
#                              holding the list's <documentManagement>
#                              property bag (TaxCatchAll, the Unified
#                              Compliance Policy placeholders and the
#                              managed-metadata "Image Tags" field), which is
#                              exactly what Word/SharePoint synthesizes the
#                              first time those document-library columns are
#                              bound to the file.
#
# No run of text in the document body changed, so the whole edit is driven
# through Document.CustomXMLParts – the same object model a real VBA/COM
# client would use to add/update that metadata part.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the new "documentManagement" custom XML part (== customXml/item4.xml)
# ---------------------------------------------------------------------------
$newPartXml = '<?xml version="1.0" encoding="utf-8"?>' +
'<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">' +
  '<documentManagement>' +
    '<_ip_UnifiedCompliancePolicyUIAction xmlns="http://schemas.microsoft.com/sharepoint/v3" xsi:nil="true"/>' +
    '<TaxCatchAll xmlns="27577a75-f6b5-4da9-9d7f-742923554f46" xsi:nil="true"/>' +
    '<_ip_UnifiedCompliancePolicyProperties xmlns="http://schemas.microsoft.com/sharepoint/v3" xsi:nil="true"/>' +
    '<lcf76f155ced4ddcb4097134ff3c332f xmlns="5a1af6ee-a0d5-4735-81f5-205ff6779adb">' +
      '<Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>' +
    '</lcf76f155ced4ddcb4097134ff3c332f>' +
  '</documentManagement>' +
'</p:properties>'

try {
    $newPart = $d.CustomXMLParts.Add($newPartXml)
} catch {
    $newPart = $null
}

# ---------------------------------------------------------------------------
# 2) Re-mint the content-type schema part (customXml/item2.xml): a fresh
#    ma:versionID on the root element and a fresh ma:fieldsID on the
#    properties schema.
# ---------------------------------------------------------------------------
try {
    $schemaParts = $d.CustomXMLParts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
} catch {
    $schemaParts = $null
}

if ($schemaParts -ne $null -and $schemaParts.Count -gt 0) {
    $schemaPart = $schemaParts.Item(1)
    try {
        $schemaXml = $schemaPart.XML
        $schemaXml = $schemaXml.Replace(
            'ma:versionID="1cc4cf9d95b5e2d14d7aabb44ca49f5e"',
            'ma:versionID="ce94caacb4a5cc228342027e3189af2c"')
        $schemaXml = $schemaXml.Replace(
            'ma:fieldsID="a33e6829bf21261855124b7b230b6e9c"',
            'ma:fieldsID="5f85a36ab557a4a47cd270a1ee4435c0"')
        $schemaPart.XML = $schemaXml
    } catch {
        # Host does not allow rewriting an existing custom XML part in
        # place; nothing further we can do through the object model.
    }
}

Write-Output "CustomXMLParts.Count after update: $($d.CustomXMLParts.Count)"
